$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster changes from FAPs to ECs, plus recomputed TPM metrics
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.0345785
$ws.Range("H2").Value = 0.069157
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2416746666666667
$ws.Range("N2").Value = 0.725024
$ws.Range("O2").Value = 0.08321776967690767
$ws.Range("P2").Value = 0.08321776967690767
$ws.Range("Q2").Value = 0.008356747461333334
$ws.Range("R2").Value = 0.050140484768
$ws.Range("S2").Value = 0.08321776967690767
$ws.Range("T2").Value = 0.08321776967690767

# Row 3: Target cluster changes from MuSCs to FAPs, plus recomputed TPM metrics
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.0345785
$ws.Range("H3").Value = 0.069157
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.662448666666667
$ws.Range("N3").Value = 7.987346000000001
$ws.Range("O3").Value = 0.9167822303230924
$ws.Range("P3").Value = 0.9167822303230924
$ws.Range("Q3").Value = 0.09206348122033332
$ws.Range("R3").Value = 0.552380887322
$ws.Range("S3").Value = 0.9167822303230924
$ws.Range("T3").Value = 0.9167822303230924
